$d = $word.ActiveDocument

function Replace-FieldWithRuns($field, $runsXml) {
    $start = $field.Code.Start - 1
    $field.Delete()
    $target = $d.Range($start, $start)
    $pkg = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" + `
        "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" + `
        "<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
        "<w:body><w:p>" + $runsXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $target.InsertXML($pkg)
}

# Field 1: "m:let self = self.name" -> "{m:" "let" " self" " " "=" " self." "name}"
$runs1 = "<w:r><w:t>{m:</w:t></w:r>" + `
         "<w:r><w:t>let</w:t></w:r>" + `
         "<w:r><w:t xml:space='preserve'> self</w:t></w:r>" + `
         "<w:r><w:t xml:space='preserve'> </w:t></w:r>" + `
         "<w:r><w:t>=</w:t></w:r>" + `
         "<w:r><w:t xml:space='preserve'> self.</w:t></w:r>" + `
         "<w:r><w:t>name}</w:t></w:r>"
Replace-FieldWithRuns $d.Fields.Item(1) $runs1

# Field 2: " m:self " -> "{" "m" ":self" "}"
$runs2 = "<w:r><w:t>{</w:t></w:r>" + `
         "<w:r><w:t>m</w:t></w:r>" + `
         "<w:r><w:t>:self</w:t></w:r>" + `
         "<w:r><w:t xml:space='preserve'>}</w:t></w:r>"
Replace-FieldWithRuns $d.Fields.Item(1) $runs2

# Field 3: " m:endlet " -> "{" "m:" "endlet" "}"
$runs3 = "<w:r><w:t>{</w:t></w:r>" + `
         "<w:r><w:t>m:</w:t></w:r>" + `
         "<w:r><w:t>endlet</w:t></w:r>" + `
         "<w:r><w:t xml:space='preserve'>}</w:t></w:r>"
Replace-FieldWithRuns $d.Fields.Item(1) $runs3

Write-Output ("Final content: " + $d.Content.Text)
